$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update - set Price (D) and Volume(1h) (E) columns
# Cells whose new value could be misread as a number are pre-formatted as
# Text so Excel keeps them as strings (matching the source data feed).

$ws.Range('D2').Value = '72.276.59'
$ws.Range('E2').Value = '  +4.46%  '
$ws.Range('D3').Value = '3.608.29'
$ws.Range('E3').Value = '  +6.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.18'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.85'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('D7').Value = '3.599.29'
$ws.Range('E7').Value = '  +6.68%  '
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  +6.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.609'
$ws.Range('E11').Value = '  +3.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.17'
$ws.Range('E12').Value = '  +3.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000292'
$ws.Range('E13').Value = '  +4.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '696.30'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').Value = '4.191.02'
$ws.Range('E15').Value = '  +6.86%  '
$ws.Range('E16').Value = '  +4.04%  '
$ws.Range('D17').Value = '72.298.39'
$ws.Range('E17').Value = '  +4.44%  '
$ws.Range('D18').Value = '3.596.11'
$ws.Range('E18').Value = '  +6.28%  '
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.51'
$ws.Range('E20').Value = '  +5.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.934'
$ws.Range('E22').Value = '  +3.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.71'
$ws.Range('E23').Value = '  +5.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.72'
$ws.Range('E24').Value = '  +3.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '105.04'
$ws.Range('E25').Value = '  +1.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.03'
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.84'
$ws.Range('E27').Value = '  +4.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('E28').Value = '  +5.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.25'
$ws.Range('E29').Value = '  +3.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.06'
$ws.Range('E30').Value = '  +4.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.52'
$ws.Range('E31').Value = '  +8.16%  '
$ws.Range('E32').Value = '  +16.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '595.87'
$ws.Range('E33').Value = '  +6.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.34'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('E35').Value = '  +1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '60.10'
$ws.Range('E36').Value = '  +2.24%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '3.670.59'
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('E39').Value = '  +5.12%  '
$ws.Range('D40').Value = '0.0₃0787'
$ws.Range('E40').Value = '  +13.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.24'
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('E42').Value = '  +6.72%  '
$ws.Range('E43').Value = '  +6.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0440'
$ws.Range('E44').Value = '  +4.30%  '
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.41'
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('E47').Value = '  +4.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.50'
$ws.Range('E48').Value = '  +6.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.133'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '133.89'
$ws.Range('E51').Value = '  +0.50%  '
